# تعديل تلقائي في شيت Card6 by admin at 2025-12-06 18:33:08
# Card6's lookup table (column A, "card" id) had rows 3-12 stuck on the
# old id "2"; they should match row 2's id "6".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

# Copy cell A2 (which already holds the correct text value "6") onto each
# of A3:A12 individually. Using Copy (rather than assigning .Value to the
# literal "6") preserves the cell's existing text data type / formatting
# instead of Excel auto-converting the numeric-looking string to a Number.
for ($r = 3; $r -le 12; $r++) {
    $ws.Range("A2").Copy($ws.Cells.Item($r, 1))
}
